$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 960.7273
$ws.Range("I8").Value = 56.8
$ws.Range("K8").Value = 170.4
$ws.Range("M8").Value = -31.39999999999998

$ws.Range("H28").Value = 1005.75
$ws.Range("I28").Value = 1013.375
$ws.Range("J28").Value = 990.5
$ws.Range("K28").Value = 1013.375
$ws.Range("L28").Value = 990.5
$ws.Range("M28").Value = -528.375
$ws.Range("N28").Value = -1960.5

$ws.Range("H31").Value = 1742
$ws.Range("I31").Value = 90.40000000000001
$ws.Range("K31").Value = 271.2
$ws.Range("M31").Value = -41.20000000000005

$ws.Range("H33").Value = 510.8889
$ws.Range("I33").Value = 157.42857
$ws.Range("K33").Value = 157.42857
$ws.Range("M33").Value = 71.57142999999999

$ws.Range("H38").Value = 747.95654
$ws.Range("I38").Value = 37.785713
$ws.Range("J38").Value = 1852.6666
$ws.Range("K38").Value = 113.357139
$ws.Range("L38").Value = 5557.9998
$ws.Range("M38").Value = 258.642861
$ws.Range("N38").Value = -6301.9998

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H87").Value = 79999.8
$ws.Range("J87").Value = 79999.8
$ws.Range("L87").Value = 79999.8
$ws.Range("N87").Value = -82495.8

$ws.Range("H90").Value = 79999.8
$ws.Range("J90").Value = 79999.8
$ws.Range("L90").Value = 239999.4
$ws.Range("N90").Value = -252479.4

$ws.Range("H95").Value = 46499.5
$ws.Range("J95").Value = 46499.5
$ws.Range("L95").Value = 46499.5
$ws.Range("N95").Value = -51991.5

$ws.Range("H107").Value = 1003
$ws.Range("I107").Value = 943
$ws.Range("K107").Value = 943
$ws.Range("M107").Value = 977

$ws.Range("H111").Value = 584.5
$ws.Range("I111").Value = 329.75
$ws.Range("K111").Value = 989.25
$ws.Range("M111").Value = 2077.75

$ws.Range("H141").Value = 1397.5
$ws.Range("I141").Value = 1397.5
$ws.Range("K141").Value = 4192.5
$ws.Range("M141").Value = 987.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3721.8
$ws.Range("I2").Value = 1203.6666
$ws.Range("K2").Value = 1203.6666
$ws.Range("M2").Value = -1090.6666

$ws.Range("H33").Value = 339333
$ws.Range("I33").Value = 339333
$ws.Range("K33").Value = 339333
$ws.Range("M33").Value = -339004

$ws.Range("H61").Value = 3469.111
$ws.Range("I61").Value = 2423.5715
$ws.Range("K61").Value = 2423.5715
$ws.Range("M61").Value = -2211.5715

$ws.Range("H74").Value = 2294.1428
$ws.Range("I74").Value = 1865
$ws.Range("J74").Value = 2866.3333
$ws.Range("K74").Value = 1865
$ws.Range("L74").Value = 2866.3333
$ws.Range("M74").Value = -991
$ws.Range("N74").Value = -4614.3333

$ws.Range("H77").Value = 2294.1428
$ws.Range("I77").Value = 1865
$ws.Range("J77").Value = 2866.3333
$ws.Range("K77").Value = 9325
$ws.Range("L77").Value = 14331.6665
$ws.Range("M77").Value = -4957
$ws.Range("N77").Value = -23067.6665

$ws.Range("H116").Value = 3721.8
$ws.Range("I116").Value = 1203.6666
$ws.Range("K116").Value = 1203.6666
$ws.Range("M116").Value = 1090.3334

$ws.Range("H136").Value = 3469.111
$ws.Range("I136").Value = 2423.5715
$ws.Range("K136").Value = 7270.7145
$ws.Range("M136").Value = -4720.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3721.8
$ws.Range("I3").Value = 1203.6666
$ws.Range("K3").Value = 1203.6666
$ws.Range("M3").Value = -1089.6666

$ws.Range("H134").Value = 2338.6667
$ws.Range("I134").Value = 2338.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7016.000100000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4481.000100000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1068.8889
$ws.Range("I58").Value = 1154.1666
$ws.Range("J58").Value = 898.3333
$ws.Range("K58").Value = 1154.1666
$ws.Range("L58").Value = 898.3333
$ws.Range("M58").Value = -951.1666
$ws.Range("N58").Value = -1304.3333

$ws.Range("H59").Value = 33999.6
$ws.Range("I59").Value = 49998
$ws.Range("K59").Value = 49998
$ws.Range("M59").Value = -48853

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H107").Value = 1327.6428
$ws.Range("I107").Value = 1394.75
$ws.Range("J107").Value = 925
$ws.Range("K107").Value = 1394.75
$ws.Range("L107").Value = 925
$ws.Range("M107").Value = 525.25
$ws.Range("N107").Value = -4765

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H134").Value = 1606.625
$ws.Range("I134").Value = 1607.5714
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 4822.7142
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -2287.7142
$ws.Range("N134").Value = -9870

$ws.Range("H136").Value = 1068.8889
$ws.Range("I136").Value = 1154.1666
$ws.Range("J136").Value = 898.3333
$ws.Range("K136").Value = 3462.4998
$ws.Range("L136").Value = 2694.9999
$ws.Range("M136").Value = -912.4998000000001
$ws.Range("N136").Value = -7794.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 10749.75
$ws.Range("I99").Value = 9333
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 27999
$ws.Range("L99").Value = 45000
$ws.Range("M99").Value = -25753
$ws.Range("N99").Value = -49492

$ws.Range("H117").Value = 2366
$ws.Range("I117").Value = 1599.75
$ws.Range("J117").Value = 3132.25
$ws.Range("K117").Value = 4799.25
$ws.Range("L117").Value = 9396.75
$ws.Range("M117").Value = -1357.25
$ws.Range("N117").Value = -16280.75

$ws.Range("H134").Value = 1212
$ws.Range("I134").Value = 1212
$ws.Range("K134").Value = 3636
$ws.Range("M134").Value = 1434

$ws.Range("H139").Value = 5561.8
$ws.Range("I139").Value = 4452.25
$ws.Range("K139").Value = 13356.75
$ws.Range("M139").Value = -8216.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H80").Value = 6750.125
$ws.Range("J80").Value = 7214.4287
$ws.Range("L80").Value = 7214.4287
$ws.Range("N80").Value = -9210.4287

$ws.Range("H83").Value = 6750.125
$ws.Range("J83").Value = 7214.4287
$ws.Range("L83").Value = 36072.14350000001
$ws.Range("N83").Value = -46056.14350000001

$ws.Range("H97").Value = 3760.4707
$ws.Range("I97").Value = 3747.25
$ws.Range("K97").Value = 3747.25
$ws.Range("M97").Value = -3251.25

$ws.Range("H101").Value = 28663.334
$ws.Range("J101").Value = 28663.334
$ws.Range("L101").Value = 28663.334
$ws.Range("N101").Value = -35153.334

$ws.Range("H122").Value = 26653.715
$ws.Range("I122").Value = 33477.727
$ws.Range("K122").Value = 100433.181
$ws.Range("M122").Value = -97983.181

$ws.Range("H132").Value = 3033.85
$ws.Range("I132").Value = 2985.611
$ws.Range("K132").Value = 8956.832999999999
$ws.Range("M132").Value = -6426.832999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1949.6666
$ws.Range("I40").Value = 1949.6666
$ws.Range("K40").Value = 1949.6666
$ws.Range("M40").Value = -1813.6666

$ws.Range("H132").Value = 2299.4
$ws.Range("J132").Value = 2499
$ws.Range("L132").Value = 7497
$ws.Range("N132").Value = -12557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H132").Value = 1786.3334
$ws.Range("I132").Value = 1771.0714
$ws.Range("K132").Value = 5313.2142
$ws.Range("M132").Value = -2783.2142

$ws.Range("H136").Value = 1134.6666
$ws.Range("I136").Value = 589.5
$ws.Range("K136").Value = 1768.5
$ws.Range("M136").Value = 781.5
